# Fixed references to new AltitudeHold controller
#
# The "Loiter" branch of the PID-controller diagram is renamed to
# "AltitudeHold": the label textbox widens (and shifts very slightly) to fit
# the new, longer text, and the connector line feeding into it from the
# summing-junction box is re-routed/resized to keep its endpoint attached.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Label textbox ("CustomShape 3") — "Loiter" -> "AltitudeHold" ----------
$label = $s.Shapes.Item("CustomShape 3")

# Reposition/resize the label box. Top/Height are unchanged (1508760 / 314640
# EMU). Left moves from 4114800 to 4114799 EMU and Width grows from 651240 to
# 1158949 EMU to fit the longer caption.
#
# Shape.Left/.Top/.Width/.Height round-trip through a single-precision float
# in this COM layer, so a handful of literal points values below are tuned
# (rather than a naive EMU/12700 division) to land exactly on the target EMU
# value once PowerPoint converts back from points to EMU on save.
$label.Left = 323.99995424999986   # -> 4114799 EMU (was 4114800)
$label.Width = 91.25586274999984   # -> 1158949 EMU (was 651240)

$label.TextFrame.TextRange.Text = "AltitudeHold"

# --- Connector ("Line 5") from the summing junction to the label -----------
$connector = $s.Shapes.Item("Line 5")

# Top is unchanged (771840 EMU); Left/Width/Height grow to reach the label's
# new position/size (3068640->3068820, 1371960->1625454, 737280->736920).
$connector.Left = 241.63938324999978    # -> 3068820 EMU (was 3068640)
$connector.Width = 127.98851324999994   # -> 1625454 EMU (was 1371960)
$connector.Height = 58.025208999999926  # -> 736920 EMU (was 737280)
